$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 24845848
$ws.Range("I112").Value = 966.6667
$ws.Range("J112").Value = 28572580
$ws.Range("K112").Value = 2900.0001
$ws.Range("L112").Value = 85717740
$ws.Range("M112").Value = -1792.0001
$ws.Range("N112").Value = -85719956
$ws.Range("H137").Value = 1562.8334
$ws.Range("I137").Value = 1427.25
$ws.Range("J137").Value = 1752.65
$ws.Range("K137").Value = 4281.75
$ws.Range("L137").Value = 5257.950000000001
$ws.Range("M137").Value = -1731.75
$ws.Range("N137").Value = -10357.95

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3420.1667
$ws.Range("I26").Value = 2504.2
$ws.Range("K26").Value = 2504.2
$ws.Range("M26").Value = -2174.2
$ws.Range("H61").Value = 4145.946
$ws.Range("I61").Value = 4738.6895
$ws.Range("J61").Value = 1997.25
$ws.Range("K61").Value = 4738.6895
$ws.Range("L61").Value = 1997.25
$ws.Range("M61").Value = -4526.6895
$ws.Range("N61").Value = -2421.25
$ws.Range("H74").Value = 1049.6863
$ws.Range("I74").Value = 870.5714
$ws.Range("J74").Value = 1885.5555
$ws.Range("K74").Value = 870.5714
$ws.Range("L74").Value = 1885.5555
$ws.Range("M74").Value = 3.42859999999996
$ws.Range("N74").Value = -3633.5555
$ws.Range("H77").Value = 1049.6863
$ws.Range("I77").Value = 870.5714
$ws.Range("J77").Value = 1885.5555
$ws.Range("K77").Value = 4352.857
$ws.Range("L77").Value = 9427.7775
$ws.Range("M77").Value = 15.14300000000003
$ws.Range("N77").Value = -18163.7775
$ws.Range("H122").Value = 3664135.8
$ws.Range("I122").Value = 5129247
$ws.Range("K122").Value = 15387741
$ws.Range("M122").Value = -15385291
$ws.Range("H136").Value = 4145.946
$ws.Range("I136").Value = 4738.6895
$ws.Range("J136").Value = 1997.25
$ws.Range("K136").Value = 14216.0685
$ws.Range("L136").Value = 5991.75
$ws.Range("M136").Value = -11666.0685
$ws.Range("N136").Value = -11091.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 38462470
$ws.Range("I107").Value = 58824332
$ws.Range("J107").Value = 1183.2222
$ws.Range("K107").Value = 58824332
$ws.Range("L107").Value = 1183.2222
$ws.Range("M107").Value = -58822412
$ws.Range("N107").Value = -5023.2222
$ws.Range("H134").Value = 3632.6897
$ws.Range("I134").Value = 4184.914
$ws.Range("J134").Value = 2792.348
$ws.Range("K134").Value = 12554.742
$ws.Range("L134").Value = 8377.044
$ws.Range("M134").Value = -10019.742
$ws.Range("N134").Value = -13447.044

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 226915.39
$ws.Range("I31").Value = 1465.3704
$ws.Range("J31").Value = 1038535.44
$ws.Range("K31").Value = 1465.3704
$ws.Range("L31").Value = 1038535.44
$ws.Range("M31").Value = -1170.3704
$ws.Range("N31").Value = -1039125.44
$ws.Range("H34").Value = 226915.39
$ws.Range("I34").Value = 1465.3704
$ws.Range("J34").Value = 1038535.44
$ws.Range("K34").Value = 1465.3704
$ws.Range("L34").Value = 1038535.44
$ws.Range("M34").Value = -1263.3704
$ws.Range("N34").Value = -1038939.44
$ws.Range("H58").Value = 1293.8167
$ws.Range("I58").Value = 803.8823
$ws.Range("J58").Value = 1934.5
$ws.Range("K58").Value = 803.8823
$ws.Range("L58").Value = 1934.5
$ws.Range("M58").Value = -600.8823
$ws.Range("N58").Value = -2340.5
$ws.Range("H132").Value = 2218.9167
$ws.Range("I132").Value = 1787.7222
$ws.Range("J132").Value = 3512.5
$ws.Range("K132").Value = 5363.1666
$ws.Range("L132").Value = 10537.5
$ws.Range("M132").Value = -2833.1666
$ws.Range("N132").Value = -15597.5
$ws.Range("H134").Value = 1884.9828
$ws.Range("I134").Value = 2121.325
$ws.Range("J134").Value = 1359.7778
$ws.Range("K134").Value = 6363.974999999999
$ws.Range("L134").Value = 4079.3334
$ws.Range("M134").Value = -3828.974999999999
$ws.Range("N134").Value = -9149.3334
$ws.Range("H136").Value = 1293.8167
$ws.Range("I136").Value = 803.8823
$ws.Range("J136").Value = 1934.5
$ws.Range("K136").Value = 2411.6469
$ws.Range("L136").Value = 5803.5
$ws.Range("M136").Value = 138.3531000000003
$ws.Range("N136").Value = -10903.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 4500
$ws.Range("I123").Value = 2000
$ws.Range("J123").Value = 5750
$ws.Range("K123").Value = 6000
$ws.Range("L123").Value = 17250
$ws.Range("M123").Value = -3550
$ws.Range("N123").Value = -22150
$ws.Range("H129").Value = 27779284
$ws.Range("I129").Value = 66667492
$ws.Range("J129").Value = 1995.1428
$ws.Range("K129").Value = 200002476
$ws.Range("L129").Value = 5985.428400000001
$ws.Range("M129").Value = -199997476
$ws.Range("N129").Value = -15985.4284
$ws.Range("H130").Value = 6765
$ws.Range("I130").Value = 530
$ws.Range("J130").Value = 13000
$ws.Range("K130").Value = 1590
$ws.Range("L130").Value = 39000
$ws.Range("M130").Value = 3430
$ws.Range("N130").Value = -49040
$ws.Range("H131").Value = 1961719.6
$ws.Range("J131").Value = 1013.26666
$ws.Range("L131").Value = 3039.79998
$ws.Range("N131").Value = -13119.79998
$ws.Range("H133").Value = 47342.44
$ws.Range("I133").Value = 130382.625
$ws.Range("J133").Value = 8264.706
$ws.Range("K133").Value = 391147.875
$ws.Range("L133").Value = 24794.118
$ws.Range("M133").Value = -386087.875
$ws.Range("N133").Value = -34914.118
$ws.Range("H134").Value = 9388.700000000001
$ws.Range("I134").Value = 11013.454
$ws.Range("K134").Value = 33040.362
$ws.Range("M134").Value = -27970.362
$ws.Range("H136").Value = 8133.0435
$ws.Range("I136").Value = 14945
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 44835
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -39735
$ws.Range("N136").Value = -23700
$ws.Range("H137").Value = 18068.334
$ws.Range("I137").Value = 8703.333000000001
$ws.Range("J137").Value = 27433.334
$ws.Range("K137").Value = 26109.999
$ws.Range("L137").Value = 82300.00199999999
$ws.Range("M137").Value = -21009.999
$ws.Range("N137").Value = -92500.00199999999
$ws.Range("H138").Value = 8759.412
$ws.Range("I138").Value = 9523.134
$ws.Range("K138").Value = 28569.402
$ws.Range("M138").Value = -23429.402
$ws.Range("H139").Value = 4053.739
$ws.Range("I139").Value = 5239.087
$ws.Range("J139").Value = 2868.3914
$ws.Range("K139").Value = 15717.261
$ws.Range("L139").Value = 8605.174199999999
$ws.Range("M139").Value = -10577.261
$ws.Range("N139").Value = -18885.1742
$ws.Range("H140").Value = 1620.8485
$ws.Range("I140").Value = 1620.8485
$ws.Range("K140").Value = 4862.5455
$ws.Range("M140").Value = 317.4544999999998
$ws.Range("H141").Value = 10466.5
$ws.Range("I141").Value = 11052.1
$ws.Range("J141").Value = 9490.5
$ws.Range("K141").Value = 33156.3
$ws.Range("L141").Value = 28471.5
$ws.Range("M141").Value = -27976.3
$ws.Range("N141").Value = -38831.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 39439536
$ws.Range("I122").Value = 70989300
$ws.Range("J122").Value = 2343.5833
$ws.Range("K122").Value = 212967900
$ws.Range("L122").Value = 7030.749899999999
$ws.Range("M122").Value = -212965450
$ws.Range("N122").Value = -11930.7499

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3543494.8
$ws.Range("I122").Value = 4204898.5
$ws.Range("J122").Value = 1669516.6
$ws.Range("K122").Value = 12614695.5
$ws.Range("L122").Value = 5008549.800000001
$ws.Range("M122").Value = -12612245.5
$ws.Range("N122").Value = -5013449.800000001
$ws.Range("H132").Value = 9545393
$ws.Range("I132").Value = 11135573
$ws.Range("J132").Value = 4308.5
$ws.Range("K132").Value = 33406719
$ws.Range("L132").Value = 12925.5
$ws.Range("M132").Value = -33404189
$ws.Range("N132").Value = -17985.5
$ws.Range("H136").Value = 7119.94
$ws.Range("I136").Value = 4971.5
$ws.Range("J136").Value = 12644.5
$ws.Range("K136").Value = 14914.5
$ws.Range("L136").Value = 37933.5
$ws.Range("M136").Value = -12364.5
$ws.Range("N136").Value = -43033.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2606.1292
$ws.Range("I122").Value = 2524.5833
$ws.Range("K122").Value = 7573.749899999999
$ws.Range("M122").Value = -5123.749899999999
$ws.Range("H123").Value = 28645.8
$ws.Range("J123").Value = 28645.8
$ws.Range("L123").Value = 28645.8
$ws.Range("N123").Value = -38445.8
$ws.Range("H132").Value = 16747
$ws.Range("I132").Value = 20256.432
$ws.Range("J132").Value = 1831.9166
$ws.Range("K132").Value = 60769.296
$ws.Range("L132").Value = 5495.7498
$ws.Range("M132").Value = -58239.296
$ws.Range("N132").Value = -10555.7498
$ws.Range("H136").Value = 9436905
$ws.Range("I136").Value = 3868.4814
$ws.Range("J136").Value = 19232752
$ws.Range("K136").Value = 11605.4442
$ws.Range("L136").Value = 57698256
$ws.Range("M136").Value = -9055.4442
$ws.Range("N136").Value = -57703356
